# Regenerate the "K" column (column G) values for rows 2-62 on Sheet1.
# The workbook stores simulated/derived counting-stat data; this edit
# overwrites the previously-generated column G ("K", formerly "Strike#")
# values with the newly regenerated ones, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 3
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 0
    41 = 2
    42 = 1
    43 = 3
    44 = 0
    45 = 0
    46 = 1
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 2
    52 = 1
    53 = 1
    54 = 1
    55 = 2
    56 = 2
    57 = 2
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
